$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add four new account test cases: acc101 - acc104 (rows 103-106)

# Row 103: acc101
$ws.Range("A103").Value = "acc101"
$ws.Range("D103").Value = "$"
$ws.Range("E103").Value = "$"
$ws.Range("F103").Value = "%"
$ws.Range("K103").Value = "%"
$ws.Range("M103").Value = "complete"

# Row 104: acc102
$ws.Range("A104").Value = "acc102"
$ws.Range("K104").Value = "%"
$ws.Range("M104").Value = "complete"

# Row 105: acc103
$ws.Range("A105").Value = "acc103"
$ws.Range("B105").Value = "$"
$ws.Range("D105").Value = "$"
$ws.Range("E105").Value = "$"
$ws.Range("F105").Value = "%"
$ws.Range("K105").Value = "%"
$ws.Range("M105").Value = "complete"

# Row 106: acc104
$ws.Range("A106").Value = "acc104"
$ws.Range("B106").Value = "$"
$ws.Range("K106").Value = "%"
$ws.Range("M106").Value = "complete"
